$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2020" column (Q), mirroring the formatting
# of the existing "2019" column (P) for every data row.
$ws.Range("P3:P34").Copy($ws.Range("Q3:Q34")) | Out-Null

# Header (year) cell.
$ws.Range("Q4").Value = 2020

# Data values for 2020, row by row (numbers where data exists, "-" where
# the source uses the placeholder dash string already present elsewhere
# in the sheet).
$ws.Range("Q5").Value = 51
$ws.Range("Q6").Value = 29
$ws.Range("Q7").Value = 22
$ws.Range("Q8").Value = 5
$ws.Range("Q9").Value = 3
$ws.Range("Q10").Value = 2
$ws.Range("Q11").Value = 15
$ws.Range("Q12").Value = 9
$ws.Range("Q13").Value = 5
$ws.Range("Q14").Value = "-"
$ws.Range("Q15").Value = "-"
$ws.Range("Q16").Value = "-"
$ws.Range("Q17").Value = "-"
$ws.Range("Q18").Value = "-"
$ws.Range("Q19").Value = "-"
$ws.Range("Q20").Value = 7
$ws.Range("Q21").Value = 7
$ws.Range("Q22").Value = "-"
$ws.Range("Q23").Value = "-"
$ws.Range("Q24").Value = "-"
$ws.Range("Q25").Value = "-"
$ws.Range("Q26").Value = 24
$ws.Range("Q27").Value = 10
$ws.Range("Q28").Value = 14
$ws.Range("Q29").Value = "-"
$ws.Range("Q30").Value = "-"
$ws.Range("Q31").Value = "-"
$ws.Range("Q32").Value = "-"
$ws.Range("Q33").Value = "-"
$ws.Range("Q34").Value = "-"

# Move the active selection to match the author's saved cursor position.
$ws.Activate()
$ws.Range("H26").Select() | Out-Null
